$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW35-FE-LIFTER ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A98").Value = 45773.45194630787
$ws1.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B98").Value = "0x01,0x90"
$ws1.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Range("D98").Value = "0x01,0x56"
$ws1.Range("E98").Value = "0xd"
$ws1.Range("F98").Value = 400
$ws1.Range("G98").Value = [double]"5.68631262647114e+23"
$ws1.Range("H98").Value = 342
$ws1.Range("I98").Value = 13

# --- Sheet 2: ROW35-MID-LIFTER ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A98").Value = 45773.30876894676
$ws2.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B98").Value = "0x01,0x90"
$ws2.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Range("D98").Value = "0x01,0x56"
$ws2.Range("E98").Value = "0xe"
$ws2.Range("F98").Value = 400
$ws2.Range("G98").Value = [double]"5.68631262647114e+23"
$ws2.Range("H98").Value = 342
$ws2.Range("I98").Value = 14

# --- Sheet 3: ROW02-FE-LIFTER ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A98").Value = 45773.4515490625
$ws3.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B98").Value = "0x01,0x90"
$ws3.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D98").Value = "0x01,0x56"
$ws3.Range("E98").Value = "0x3"
$ws3.Range("F98").Value = 400
$ws3.Range("G98").Value = [double]"5.68631262647114e+23"
$ws3.Range("H98").Value = 342
$ws3.Range("I98").Value = 3

# --- Sheet 4: ROW02-MID-LIFTER ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A98").Value = 45773.51554211805
$ws4.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B98").Value = "0x01,0x90"
$ws4.Range("C98").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D98").Value = "0x01,0x56"
$ws4.Range("E98").Value = "0x3"
$ws4.Range("F98").Value = 400
$ws4.Range("G98").Value = [double]"9.85046333984776e+23"
$ws4.Range("H98").Value = 342
$ws4.Range("I98").Value = 3
